$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new player rows (4-8), copying the existing data-row style (row 3) down first so the
# new rows inherit the same formatting as the current data rows without minting unexpected styles ---
$ws.Range("A3:F3").Copy()
$ws.Range("A4:F8").PasteSpecial(-4122)

$data = @(
  @("Anna","Geiger",1,(Get-Date -Year 1995 -Month 5 -Day 27 -Hour 0 -Minute 0 -Second 0),"Winger","Hells Teddies"),
  @("Laura","Anninger",4,(Get-Date -Year 1994 -Month 7 -Day 17 -Hour 0 -Minute 0 -Second 0),"Attacking Midfielder","Elements"),
  @("Phillipp","Stöllinger",9,(Get-Date -Year 1992 -Month 7 -Day 1 -Hour 0 -Minute 0 -Second 0),"Central Midfielder","Valantic"),
  @("Sebastian","Meier",8,(Get-Date -Year 1993 -Month 3 -Day 6 -Hour 0 -Minute 0 -Second 0),"Center Back","New Team"),
  @("Florian","Forsthuber",3,(Get-Date -Year 1991 -Month 11 -Day 2 -Hour 0 -Minute 0 -Second 0),"Sweeper","Alchimiste")
)

$r = 4
foreach ($row in $data) {
  $ws.Cells.Item($r,1).Value = $row[0]
  $ws.Cells.Item($r,2).Value = $row[1]
  $ws.Cells.Item($r,3).Value = $row[2]
  $ws.Cells.Item($r,4).Value = $row[3]
  $ws.Cells.Item($r,5).Value = $row[4]
  $ws.Cells.Item($r,6).Value = $row[5]
  $r++
}

# --- Widen the player_position column (E) ---
$ws.Columns.Item(5).ColumnWidth = 14.38

# --- Bold the header row and extend the header formatting out through column Z ---
$ws.Range("A1:F1").Font.Bold = $true
$ws.Range("A1").Copy()
$ws.Range("G1:Z1").PasteSpecial(-4122)
